$wb = $excel.ActiveWorkbook

# The "Config" sheet holds app settings as Key/Value rows. Row 2 previously
# stored a stray portfolio-id string under the "reinvest_dividends" key
# (a multi-currency import bug). Fix it so the key is correctly named
# "reinvested_dividends" and its value is the real JSON config payload.
$ws = $wb.Worksheets.Item(4)

# Set the value cell before the key cell so the shared-string table grows
# in the same order as the fixed workbook.
$ws.Range("B2").Value = '[{"portfolio_id": "9e792bb8-94e7-4ed3-b8cc-43b50d34c337", "symbol": "ACME"}]'
$ws.Range("A2").Value = "reinvested_dividends"

# Match the saved selection/active cell on the Config sheet.
$ws.Activate()
$ws.Range("B7").Select()
